# Add new SCT-rank personnel rows to the ALPHA sheet (flight personnel roster).
# Commit message: "added some statuses :)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALPHA")

# Columns: A=rank, B=displayName, C=sheetName, D=commSec, E=nor
$newRows = @(
    @("SCT", "BRAYDEN",    "BRAYDEN",    "NIL", "REGULAR"),
    @("SCT", "NOAH LAM",   "NOAH LAM",   "NIL", "NSF"),
    @("SCT", "MARCUS",     "MARCUS",     "NIL", "NSF"),
    @("SCT", "MENG LONG",  "MENG LONG",  "NIL", "NSF"),
    @("SCT", "KAI",        "KAI",        "NIL", "NSF"),
    @("SCT", "CHARLES",    "CHARLES",    "NIL", "NSF"),
    @("SCT", "DARSHAN",    "DARSHAN",    "NIL", "NSF"),
    @("SCT", "ZHONG PING", "ZHONG PING", "NIL", "NSF"),
    @("SCT", "DERRILL",    "DERRILL",    "NIL", "NSF")
)

$startRow = 33
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
}
